# Update ESW_data.xlsx per upload diff:
#  - D6: replace the stray address fragment with the RIT chapter name
#  - D8, D28, D32, D37, D41: the "Chapter_Locator" column had accidentally
#    been filled with a duplicate of the Chapter_Name column; fix these
#    five rows so Chapter_Locator holds the street address (same as
#    column B) like the rest of the sheet already does.
#  - Column B (Address) widens to fit its longest value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = "Rochester institute of technology (RIT)"

$ws.Range("D8").Value  = $ws.Range("B8").Value2
$ws.Range("D28").Value = $ws.Range("B28").Value2
$ws.Range("D32").Value = $ws.Range("B32").Value2
$ws.Range("D37").Value = $ws.Range("B37").Value2
$ws.Range("D41").Value = $ws.Range("B41").Value2

$ws.Columns.Item(2).ColumnWidth = 40.83203125
